$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.365.22'
$ws.Range("E2").Value = '  +1.27%  '

$ws.Range("D3").Value = '1.919.41'
$ws.Range("E3").Value = '  +0.57%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.0000'
$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.8116'
$ws.Range("E5").Value = '  +4.01%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.35'
$ws.Range("E6").Value = '  +1.10%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9997'
$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3242'
$ws.Range("E8").Value = '  +2.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.23'
$ws.Range("E9").Value = '  +3.93%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07205'
$ws.Range("E10").Value = '  +4.88%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7880'
$ws.Range("E11").Value = '  +6.50%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08095'
$ws.Range("E12").Value = '  +1.48%  '

$ws.Range("D13").Value = '1.931.47'
$ws.Range("E13").Value = '  +1.22%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.421'
$ws.Range("E14").Value = '  +4.37%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.87'
$ws.Range("E15").Value = '  +2.31%  '

$ws.Range("D16").Value = '30.368.01'
$ws.Range("E16").Value = '  +1.26%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.33'
$ws.Range("E17").Value = '  +2.96%  '

$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '254.43'
$ws.Range("E18").Value = '  +3.62%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.051'
$ws.Range("E19").Value = '  +3.12%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007837'
$ws.Range("E20").Value = '  +1.45%  '

$ws.Range("D21").Value = '2.174.96'
$ws.Range("E21").Value = '  +0.92%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  +0.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.997'
$ws.Range("E23").Value = '  +16.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("E25").Value = '  +17.88%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.526'
$ws.Range("E26").Value = '  +3.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '167.69'
$ws.Range("E27").Value = '  -0.49%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.16'
$ws.Range("E28").Value = '  +1.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.145'
$ws.Range("E29").Value = '  +5.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.376'
$ws.Range("E30").Value = '  +0.69%  '

$ws.Range("E31").Value = '  +1.42%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.357'
$ws.Range("E32").Value = '  +1.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.144'
$ws.Range("E33").Value = '  +1.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05621'
$ws.Range("E34").Value = '  +0.99%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.305'
$ws.Range("E35").Value = '  +4.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7439'
$ws.Range("E36").Value = '  +1.58%  '

$ws.Range("E37").Value = '  +0.17%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.715'
$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01959'
$ws.Range("E39").Value = '  +1.75%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.823'
$ws.Range("E40").Value = '  +1.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4499'
$ws.Range("E41").Value = '  +2.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.80'
$ws.Range("E42").Value = '  +2.62%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.991'
$ws.Range("E43").Value = '  -2.21%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8552'
$ws.Range("E44").Value = '  +1.57%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.938'
$ws.Range("E45").Value = '  +3.81%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9997'
$ws.Range("E46").Value = '  -0.06%  '

$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.27'
$ws.Range("E47").Value = '  +2.85%  '

$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.030.65'
$ws.Range("E48").Value = '  +4.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '10.02'
$ws.Range("E49").Value = '  +3.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.670'
$ws.Range("E50").Value = '  +2.17%  '

$ws.Range("D51").Value = '2.064.64'
$ws.Range("E51").Value = '  +0.42%  '

